$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Kaviya"
$ws.Range("B3").Value = "chennai"
$ws.Range("C3").Value = "'123"
